$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 244, shifting existing rows 244:309 down to 245:310.
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new record.
$ws.Range("A244").Value = 8
$ws.Range("B244").Value = "Terminal La Palmera de La Serena"
$ws.Range("C244").Value = "Coquimbo"
$ws.Range("D244").Value = 44876
$ws.Range("E244").Value = 4
$ws.Range("F244").Value = "Fruta"
$ws.Range("G244").Value = 100103
$ws.Range("H244").Value = "Frutos de hueso (carozo)"
$ws.Range("I244").Value = 100103001
$ws.Range("J244").Value = "Cereza"
$ws.Range("K244").Value = "Early Burlat"
$ws.Range("L244").Value = "Primera"
$ws.Range("M244").Value = 160
$ws.Range("N244").Value = 31000
$ws.Range("O244").Value = 32000
$ws.Range("P244").Value = 31500
$ws.Range("Q244").Value = "`$/caja 15 kilos"
$ws.Range("R244").Value = "Región de O'Higgins"
$ws.Range("S244").Value = 2100
$ws.Range("T244").Value = 15

# Match the date-number style used by the rest of column D.
$ws.Range("D244").NumberFormat = $ws.Range("D245").NumberFormat
